$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column C (E values) and the new A-column iteration labels (rows 7-33)
# to be treated as text so Excel does not silently convert numeric-looking
# strings (e.g. "13.2", "6", "7.20734747345283e-05") into real numbers.
$ws.Range("A7:A33").NumberFormat = "@"
$ws.Range("C2:C33").NumberFormat = "@"

$ws.Range("B2").Value = '[2.5 0.25 1;1.66666666666667 3.16666666666667 2.16666666666667;1.77083333333333 2.80208333333333 2.58333333333333]'
$ws.Range("C2").Value = '13.2'
$ws.Range("B3").Value = '[0.364583333333333 -1.27604166666667 -0.416666666666667;2.79513888888889 3.71701388888889 3.18055555555556;2.82335069444444 3.61827256944444 3.21527777777778]'
$ws.Range("C3").Value = '10.4530612244898'
$ws.Range("B4").Value = '[-1.00802951388889 -2.09689670138889 -1.49305555555556;3.53479456018518 4.12821903935185 3.79282407407407;3.49631980613426 4.02850567853009 3.7349537037037]'
$ws.Range("C4").Value = '2.47407129452994'
$ws.Range("B5").Value = '[-1.89925582320602 -2.61041711877894 -2.21209490740741;4.0167839144483 4.40219379943094 4.18557098765432;3.93283891089169 4.2809622493791 4.08616657021605]'
$ws.Range("C5").Value = '0.991019001341996'
$ws.Range("B6").Value = '[-2.47900739128207 -2.94212647426276 -2.68226152584877;4.3305317757061 4.58125727461199 4.4404799221965;4.21674667562477 4.44351472776123 4.31629347310635]'
$ws.Range("C6").Value = '0.521896549517097'
$ws.Range("A7").Value = '6'
$ws.Range("B7").Value = '[-2.85627216959196 -3.1577003198396 -2.98850667820055;4.53472366712384 4.69788109193286 4.60628887328264;4.40148918921401 4.54909242691654 4.46624445555469]'
$ws.Range("C7").Value = '0.302826485080389'
$ws.Range("A8").Value = '7'
$ws.Range("B8").Value = '[-3.10178734494989 -3.29795703240792 -3.18783888273933;4.66761003176426 4.77378928378585 4.71418517923377;4.52171458265262 4.61777582430848 4.56385300690364]'
$ws.Range("C8").Value = '0.1842101154533'
$ws.Range("A9").Value = '8'
$ws.Range("B9").Value = '[-3.2615648151495 -3.38922987499363 -3.31756538787715;4.75409077965756 4.82319061261101 4.78440142410082;4.59995531455405 4.66247101871827 4.62737801139801]'
$ws.Range("C9").Value = '0.115021163090169'
$ws.Range("A10").Value = '9'
$ws.Range("B10").Value = '[-3.3655457420202 -3.44862846881739 -3.40199007377462;4.81037127558779 4.85534047609188 4.83009704728341;4.65087326986568 4.6915576739879 4.66871953428829]'
$ws.Range("C10").Value = '0.0729358060726012'
$ws.Range("A11").Value = '10'
$ws.Range("B11").Value = '[-3.43321509162368 -3.48728419406286 -3.4569325526067;4.84699784943817 4.87626318371059 4.85983511268975;4.68400996990526 4.71048682536164 4.69562406720675]'
$ws.Range("C11").Value = '0.0466883823364352'
$ws.Range("A12").Value = '11'
$ws.Range("B12").Value = '[-3.47725337203126 -3.51244080046376 -3.49268836812069;4.8708339197033 4.88987939608224 4.87918823421267;4.70557487759371 4.72280565126929 4.71313317152226]'
$ws.Range("C12").Value = '0.0300641513603654'
$ws.Range("A13").Value = '12'
$ws.Range("B13").Value = '[-3.50591287857433 -3.52881237269632 -3.51595776142063;4.88634610611727 4.89874063991933 4.89178297902671;4.71960902257964 4.73082257295537 4.72452785613122]'
$ws.Range("C13").Value = '0.0194322435395551'
$ws.Range("A14").Value = '13'
$ws.Range("B14").Value = '[-3.52456409087777 -3.53946676641718 -3.53110116233564;4.89644122348857 4.90450741545223 4.89997946553522;4.72874225092646 4.73603987514768 4.73194336007597]'
$ws.Range("C14").Value = '0.0125905250928361'
$ws.Range("A15").Value = '14'
$ws.Range("B15").Value = '[-3.53670204307966 -3.54640049916301 -3.5409562791894;4.9030109868987 4.90826035358406 4.90531362611361;4.73468603020011 4.73943522358087 4.73676926796498]'
$ws.Range("C15").Value = '0.008170323437991'
$ws.Range("A16").Value = '15'
$ws.Range("B16").Value = '[-3.54460125527408 -3.55091287697848 -3.54736985356769;4.9072864984827 4.91070271405551 4.9087850243843;4.73855415992562 4.74164486959767 4.73990990238373]'
$ws.Range("C16").Value = '0.0053072650126202'
$ws.Range("A17").Value = '16'
$ws.Range("B17").Value = '[-3.54974195382484 -3.55384947034047 -3.55154371948009;4.91006894256229 4.91229216862737 4.91104416258944;4.74107148549995 4.74308287680595 4.7419537840277]'
$ws.Range("C17").Value = '0.0034497285738093'
$ws.Range("A18").Value = '17'
$ws.Range("B18").Value = '[-3.55308744967169 -3.5557605648735 -3.55426001395593;4.91187971886447 4.91332656378134 4.91251437863267;4.74270972632869 4.7440187121006 4.74328391406429]'
$ws.Range("C18").Value = '0.0022432749399862'
$ws.Range("A19").Value = '18'
$ws.Range("B19").Value = '[-3.5552646523127 -3.55700427888631 -3.55602774100664;4.91305814715368 4.91399973390744 4.91347117499371;4.74377587090701 4.74462774082708 4.74414954438072]'
$ws.Range("C19").Value = '0.0014591476615323'
$ws.Range("A20").Value = '19'
$ws.Range("B20").Value = '[-3.55668154581877 -3.55781367084412 -3.55717815343565;4.91382505206134 4.91443782375823 4.91409384489364;4.74446970312139 4.74502408833802 4.74471288467387]'
$ws.Range("C20").Value = '0.0009492780527309'
$ws.Range("A21").Value = '20'
$ws.Range("B21").Value = '[-3.5576036406067 -3.55834041198768 -3.55792682600717;4.91432414321757 4.91472292660212 4.91449906989247;4.7449212395748 4.74528202584177 4.74507949878136]'
$ws.Range("C21").Value = '0.0006176436879472'
$ws.Range("A22").Value = '21'
$ws.Range("B22").Value = '[-3.55820372720058 -3.55868320787248 -3.55841405181003;4.91464894487125 4.91490846760802 4.91476278474313;4.74521509328255 4.74544988801829 4.74531808619549]'
$ws.Range("C22").Value = '0.0004018975321017'
$ws.Range("A23").Value = '22'
$ws.Range("B23").Value = '[-3.55859425529471 -3.55890629471516 -3.55873113165509;4.91486032131605 4.91502921514039 4.91493440673748;4.74540632923018 4.74555913041188 4.74547335560006]'
$ws.Range("C23").Value = '0.0002615254675413'
$ws.Range("A24").Value = '23'
$ws.Range("B24").Value = '[-3.55884840560213 -3.55905147656123 -3.55893748285314;4.91499788219639 4.91510779597218 4.91504609596875;4.74553078295223 4.74563022385773 4.74557440279102]'
$ws.Range("C24").Value = '0.0001701870476013'
$ws.Range("A25").Value = '24'
$ws.Range("B25").Value = '[-3.55901380312341 -3.559145958908 -3.55907177337208;4.91508740492357 4.91515893529571 4.91511878178288;4.74561177572124 4.74567649049357 4.74564016291183]'
$ws.Range("C25").Value = '0.0001107510993776'
$ws.Range("A26").Value = '25'
$ws.Range("B26").Value = '[-3.55912144155329 -3.55920744671857 -3.55915916779307;4.91514566508199 4.91519221606345 4.91516608471008;4.74566448470031 4.74570660018324 4.74568295869315]'
$ws.Range("C26").Value = '7.20734747345283e-05'
$ws.Range("A27").Value = '26'
$ws.Range("B27").Value = '[-3.55919149116165 -3.55924746213921 -3.55921604287913;4.91518357999105 4.91521387472893 4.9151968688039;4.74569878697827 4.74572619515477 4.74571080959848]'
$ws.Range("C27").Value = '4.69036560347471e-05'
$ws.Range("A28").Value = '27'
$ws.Range("B28").Value = '[-3.55923707848242 -3.55927350362408 -3.55925305640216;4.9152082544919 4.91522796989026 4.91521690266836;4.74572111042854 4.74573894729249 4.74572893458426]'
$ws.Range("C28").Value = '3.05239279149347e-05'
$ws.Range("A29").Value = '28'
$ws.Range("B29").Value = '[-3.55926674608319 -3.55929045106394 -3.5592771442934;4.91522431231737 4.91523714282721 4.91522994043156;4.74573563822265 4.74574724620816 4.74574073007549]'
$ws.Range("C29").Value = '1.98644141255277e-05'
$ws.Range("A30").Value = '29'
$ws.Range("B30").Value = '[-3.55928605334936 -3.55930148022449 -3.55929282036141;4.91523476252913 4.9152431124483 4.91523842522836;4.74574509271107 4.74575264702823 4.74574840641879]'
$ws.Range("C30").Value = '1.29274288221886e-05'
$ws.Range("A31").Value = '30'
$ws.Range("B31").Value = '[-3.55929861825238 -3.55930865785034 -3.55930302213067;4.91524156338308 4.91524699739552 4.91524394701731;4.74575124556197 4.74575616180758 4.74575340207734]'
$ws.Range("C31").Value = '8.41296796716408e-06'
$ws.Range("A32").Value = '31'
$ws.Range("B32").Value = '[-3.55930679531829 -3.55931332895043 -3.55930966130165;4.9152459892852 4.91524952566569 4.91524754052154;4.74575524975263 4.7457584491776 4.74575665318315]'
$ws.Range("C32").Value = '5.47503325497692e-06'
$ws.Range("A33").Value = '32'
$ws.Range("B33").Value = '[-3.55931211684021 -3.55931636883807 -3.55931398198273;4.91524886960137 4.91525117102911 4.91524987912463;4.74575785562479 4.74575993776651 4.74575876895805]'
$ws.Range("C33").Value = '3.56307183135093e-06'

# Reset the number format back to the default style so the cells end up
# with no explicit style index, matching a plain inline/shared string cell.
$ws.Range("A7:A33").Style = "Normal"
$ws.Range("C2:C33").Style = "Normal"

Write-Output "done"
